$d = $word.ActiveDocument

$replacements = @(
    @("636÷2=", "996÷7="),
    @("946÷3=", "794÷8="),
    @("265÷3=", "524÷8="),
    @("202÷3=", "806÷3="),
    @("709÷6=", "330÷9="),
    @("374÷7=", "645÷9="),
    @("158÷7=", "193÷9="),
    @("144÷5=", "540÷9="),
    @("927÷9=", "672÷9="),
    @("890÷8=", "173÷4="),
    @("527÷6=", "958÷5="),
    @("246÷8=", "338÷4="),
    @("228÷2=", "360÷5="),
    @("689÷4=", "200÷2="),
    @("471÷6=", "384÷2="),
    @("574÷2=", "552÷5="),
    @("338÷5=", "122÷3="),
    @("257÷4=", "111÷5="),
    @("307÷7=", "801÷4="),
    @("512÷6=", "404÷5="),
    @("512÷7=", "187÷7="),
    @("330÷4=", "489÷5="),
    @("461÷2=", "461÷5="),
    @("794÷6=", "694÷3="),
    @("540÷4=", "538÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
